# Update countries & provincias Spain
# - Refresh a handful of per-country COVID case counters (India, Honduras,
#   Bielorrusia, Tailandia, Mongolia, Butan rows).
# - Honduras overtakes Bielorrusia in total cases, so the two swap places
#   in the "Casos totales" (column B) descending sort.
# - Bump the "Datos actualizados" timestamp in the title cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: India ---------------------------------------------------------
$ws.Range("B5").Value = 7120538
$ws.Range("C5").Value = 1238
$ws.Range("D5").Value = 6149535
$ws.Range("E5").Value = 861819

# --- Rows 54/55: Honduras & Bielorrusia swap order ------------------------
# Honduras' case count now exceeds Bielorrusia's, so Honduras moves up to
# row 54 (with its freshly updated figures) and Bielorrusia drops to row 55
# keeping the figures Honduras' old row used to hold.
$ws.Range("A54").Value = "Honduras"
$ws.Range("B54").Value = 84081
$ws.Range("C54").Value = 935
$ws.Range("D54").Value = 32012
$ws.Range("E54").Value = 49557
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 8
$ws.Range("H54").Value = 2512

$ws.Range("A55").Value = "Bielorrusia"
$ws.Range("B55").Value = 83534
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 77220
$ws.Range("E55").Value = 5418
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 896

# --- Row 143: Tailandia ----------------------------------------------------
$ws.Range("B143").Value = 3641
$ws.Range("C143").Value = 5
$ws.Range("D143").Value = 3454
$ws.Range("E143").Value = 128

# --- Row 186: Mongolia ------------------------------------------------------
$ws.Range("B186").Value = 318
$ws.Range("C186").Value = 2
$ws.Range("E186").Value = 8

# --- Row 187: Butan ----------------------------------------------------------
$ws.Range("B187").Value = 309
$ws.Range("C187").Value = 3
$ws.Range("D187").Value = 288
$ws.Range("E187").Value = 21

# --- Title timestamp ---------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 12 de Octubre de 2020 a las 06:35"
